# Travail sur les religions
# Sort the "Tableau1" table on the worksheet by column A (Id), ascending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:K11")
$sortRange.Sort($ws.Range("A1:A11"), 1, $null, $null, 1, $null, 1, 1)

$ws.Range("B11").Select()
